$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = 13.03
$ws.Range("E6").Value = 12.487
$ws.Range("D7").Value = -7.251
$ws.Range("E7").Value = 12.608
$ws.Range("A8").Value = -21.107
$ws.Range("E8").Value = 12.766
$ws.Range("E9").Value = 12.3
$ws.Range("A10").Value = -20.712
$ws.Range("E10").Value = 12.54
$ws.Range("A12").Value = -21.766
$ws.Range("E12").Value = 13.055
$ws.Range("B13").Value = 6.448
$ws.Range("A18").Value = -21.766
$ws.Range("D20").Value = -8.222
$ws.Range("A25").Value = -21.753
